# Update functions and Data Model (#50)
# Adds a new "Authorship Resource" column to Table1 / Sheet1, fills it in for
# every data row, tidies up the "Footnote" column width, and lets row heights
# re-settle for the few rows whose content no longer drives the tallest cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Add the new table column (becomes column O, id 15) -----------------
$newCol = $lo.ListColumns.Add()
$ws.Range("O1").Value = "Authorship Resource"

# --- 2. Fill in the authorship info for every data row ----------------------
$authors = "Noémi Villars-Amberg, Daniela Subotic"
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 15).Value = $authors
}

# --- 3. Match the data-cell formatting used by the other text columns ------
#        (Arial 14, "@" text format, top-aligned + wrapped) by copying the
#        format already applied on the neighbouring "Alternative Name" column.
$ws.Range("L2:L55").Copy()
$ws.Range("O2:O55").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Narrow the "Footnote" column (G) now that it's visually separate
#        from the wide "Description" column (F) ------------------------------
$ws.Columns.Item(7).ColumnWidth = 18.5

# --- 5. Let row heights re-settle: a handful of short rows are now driven by
#        the new column's wrapped text instead of their old tallest cell -----
$ws.Rows.Item(42).RowHeight = 95
$ws.Rows.Item(43).RowHeight = 95
$ws.Rows.Item(50).RowHeight = 95
$ws.Rows.Item(53).RowHeight = 95

# --- 6. Leave the selection on the freshly-filled column --------------------
$ws.Range("O2:O55").Select()
